$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 852
$ws.Cells.Item(18, 9).Value = 852
$ws.Cells.Item(18, 11).Value = 852
$ws.Cells.Item(18, 13).Value = -568
$ws.Cells.Item(20, 8).Value = 4943.6665
$ws.Cells.Item(20, 10).Value = 12631
$ws.Cells.Item(20, 12).Value = 12631
$ws.Cells.Item(20, 14).Value = -13091
$ws.Cells.Item(33, 8).Value = 78.888885
$ws.Cells.Item(33, 9).Value = 81.25
$ws.Cells.Item(33, 11).Value = 81.25
$ws.Cells.Item(33, 13).Value = 147.75
$ws.Cells.Item(35, 8).Value = 4943.6665
$ws.Cells.Item(35, 10).Value = 12631
$ws.Cells.Item(35, 12).Value = 12631
$ws.Cells.Item(35, 14).Value = -13389
$ws.Cells.Item(40, 8).Value = 6027.25
$ws.Cells.Item(40, 10).Value = 8666
$ws.Cells.Item(40, 12).Value = 8666
$ws.Cells.Item(40, 14).Value = -9016
$ws.Cells.Item(88, 8).Value = 1644.25
$ws.Cells.Item(88, 9).Value = 868.3333
$ws.Cells.Item(88, 10).Value = 2109.8
$ws.Cells.Item(88, 11).Value = 868.3333
$ws.Cells.Item(88, 12).Value = 2109.8
$ws.Cells.Item(88, 13).Value = -462.3333
$ws.Cells.Item(88, 14).Value = -2921.8
$ws.Cells.Item(91, 8).Value = 1644.25
$ws.Cells.Item(91, 9).Value = 868.3333
$ws.Cells.Item(91, 10).Value = 2109.8
$ws.Cells.Item(91, 11).Value = 868.3333
$ws.Cells.Item(91, 12).Value = 2109.8
$ws.Cells.Item(91, 13).Value = 535.6667
$ws.Cells.Item(91, 14).Value = -4917.8
$ws.Cells.Item(101, 8).Value = 674.125
$ws.Cells.Item(101, 10).Value = 542.3333
$ws.Cells.Item(101, 12).Value = 1626.9999
$ws.Cells.Item(101, 14).Value = -4870.9999
$ws.Cells.Item(105, 8).Value = 31223.334
$ws.Cells.Item(105, 10).Value = 31223.334
$ws.Cells.Item(105, 12).Value = 31223.334
$ws.Cells.Item(105, 14).Value = -38211.334
$ws.Cells.Item(121, 8).Value = 760.5714
$ws.Cells.Item(121, 10).Value = 760.5714
$ws.Cells.Item(121, 12).Value = 2281.7142
$ws.Cells.Item(121, 14).Value = -5775.7142
$ws.Cells.Item(129, 8).Value = 1827.4
$ws.Cells.Item(129, 9).Value = 1182.875
$ws.Cells.Item(129, 10).Value = 2564
$ws.Cells.Item(129, 11).Value = 3548.625
$ws.Cells.Item(129, 12).Value = 7692
$ws.Cells.Item(129, 13).Value = 1451.375
$ws.Cells.Item(129, 14).Value = -17692
$ws.Cells.Item(135, 8).Value = 2529.375
$ws.Cells.Item(135, 9).Value = 2580
$ws.Cells.Item(135, 11).Value = 23220
$ws.Cells.Item(135, 13).Value = -20685
$ws.Cells.Item(138, 8).Value = 2317.7856
$ws.Cells.Item(138, 10).Value = 2951.25
$ws.Cells.Item(138, 12).Value = 8853.75
$ws.Cells.Item(138, 14).Value = -19133.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 601.6667
$ws.Cells.Item(22, 9).Value = 601.6667
$ws.Cells.Item(22, 11).Value = 601.6667
$ws.Cells.Item(22, 13).Value = -302.6667
$ws.Cells.Item(45, 8).Value = 3147.5454
$ws.Cells.Item(45, 9).Value = 2328
$ws.Cells.Item(45, 11).Value = 2328
$ws.Cells.Item(45, 13).Value = -1951
$ws.Cells.Item(46, 8).Value = 3374.6
$ws.Cells.Item(46, 9).Value = 2791.3333
$ws.Cells.Item(46, 10).Value = 4249.5
$ws.Cells.Item(46, 11).Value = 2791.3333
$ws.Cells.Item(46, 12).Value = 4249.5
$ws.Cells.Item(46, 13).Value = -2472.3333
$ws.Cells.Item(46, 14).Value = -4887.5
$ws.Cells.Item(97, 8).Value = 1531.3077
$ws.Cells.Item(97, 9).Value = 1079.5
$ws.Cells.Item(97, 11).Value = 1079.5
$ws.Cells.Item(97, 13).Value = -583.5
$ws.Cells.Item(110, 8).Value = 3045.3
$ws.Cells.Item(110, 9).Value = 2315.7144
$ws.Cells.Item(110, 10).Value = 4747.6665
$ws.Cells.Item(110, 11).Value = 2315.7144
$ws.Cells.Item(110, 12).Value = 4747.6665
$ws.Cells.Item(110, 13).Value = -270.7143999999998
$ws.Cells.Item(110, 14).Value = -8837.666499999999
$ws.Cells.Item(122, 8).Value = 1259.5
$ws.Cells.Item(122, 9).Value = 1146.7273
$ws.Cells.Item(122, 11).Value = 3440.1819
$ws.Cells.Item(122, 13).Value = -990.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(31, 8).Value = 95
$ws.Cells.Item(31, 9).Value = 95
$ws.Cells.Item(31, 11).Value = 95
$ws.Cells.Item(31, 13).Value = 157
$ws.Cells.Item(62, 8).Value = 100000
$ws.Cells.Item(62, 10).Value = 100000
$ws.Cells.Item(62, 12).Value = 100000
$ws.Cells.Item(62, 14).Value = -101372
$ws.Cells.Item(65, 8).Value = 100000
$ws.Cells.Item(65, 10).Value = 100000
$ws.Cells.Item(65, 12).Value = 300000
$ws.Cells.Item(65, 14).Value = -306864
$ws.Cells.Item(94, 8).Value = 619.2308
$ws.Cells.Item(94, 9).Value = 614.2727
$ws.Cells.Item(94, 11).Value = 614.2727
$ws.Cells.Item(94, 13).Value = -163.2727
$ws.Cells.Item(107, 8).Value = 3612
$ws.Cells.Item(107, 9).Value = 1307.0625
$ws.Cells.Item(107, 11).Value = 1307.0625
$ws.Cells.Item(107, 13).Value = 612.9375
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4475.5264
$ws.Cells.Item(31, 9).Value = 1813.35
$ws.Cells.Item(31, 10).Value = 7433.5
$ws.Cells.Item(31, 11).Value = 1813.35
$ws.Cells.Item(31, 12).Value = 7433.5
$ws.Cells.Item(31, 13).Value = -1518.35
$ws.Cells.Item(31, 14).Value = -8023.5
$ws.Cells.Item(34, 8).Value = 4475.5264
$ws.Cells.Item(34, 9).Value = 1813.35
$ws.Cells.Item(34, 10).Value = 7433.5
$ws.Cells.Item(34, 11).Value = 1813.35
$ws.Cells.Item(34, 12).Value = 7433.5
$ws.Cells.Item(34, 13).Value = -1611.35
$ws.Cells.Item(34, 14).Value = -7837.5
$ws.Cells.Item(105, 8).Value = 1633.1578
$ws.Cells.Item(105, 9).Value = 1633.1578
$ws.Cells.Item(105, 11).Value = 1633.1578
$ws.Cells.Item(105, 13).Value = 113.8422

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(94, 8).Value = 29500
$ws.Cells.Item(94, 10).Value = 29500
$ws.Cells.Item(94, 12).Value = 88500
$ws.Cells.Item(94, 14).Value = -89852
$ws.Cells.Item(107, 8).Value = 985.6667
$ws.Cells.Item(107, 10).Value = 942.0714
$ws.Cells.Item(107, 12).Value = 2826.2142
$ws.Cells.Item(107, 14).Value = -6666.2142
$ws.Cells.Item(121, 8).Value = 1854.6364
$ws.Cells.Item(121, 9).Value = 480.8
$ws.Cells.Item(121, 11).Value = 1442.4
$ws.Cells.Item(121, 13).Value = -132.4000000000001
$ws.Cells.Item(131, 8).Value = 2971.4285
$ws.Cells.Item(131, 9).Value = 2800
$ws.Cells.Item(131, 10).Value = 3000
$ws.Cells.Item(131, 11).Value = 8400
$ws.Cells.Item(131, 12).Value = 9000
$ws.Cells.Item(131, 13).Value = -3360
$ws.Cells.Item(131, 14).Value = -19080

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 3666.6667
$ws.Cells.Item(70, 8).Value = 2500
$ws.Cells.Item(70, 9).Value = 2500
$ws.Cells.Item(70, 11).Value = 2500
$ws.Cells.Item(70, 13).Value = -2230
$ws.Cells.Item(73, 8).Value = 2500
$ws.Cells.Item(73, 9).Value = 2500
$ws.Cells.Item(73, 11).Value = 2500
$ws.Cells.Item(73, 13).Value = -1564
$ws.Cells.Item(80, 8).Value = 2250
$ws.Cells.Item(80, 9).Value = 2250
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 2250
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = -1252
$ws.Cells.Item(80, 14).ClearContents()
$ws.Cells.Item(83, 8).Value = 2250
$ws.Cells.Item(83, 9).Value = 2250
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 11250
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = -6258
$ws.Cells.Item(83, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 5847.4
$ws.Cells.Item(113, 9).Value = 3249.8333
$ws.Cells.Item(113, 10).Value = 7579.1113
$ws.Cells.Item(113, 11).Value = 3249.8333
$ws.Cells.Item(113, 12).Value = 7579.1113
$ws.Cells.Item(113, 13).Value = -1079.8333
$ws.Cells.Item(113, 14).Value = -11919.1113
$ws.Cells.Item(122, 8).Value = 175349.31
$ws.Cells.Item(122, 10).Value = 3286.5454
$ws.Cells.Item(122, 12).Value = 9859.636200000001
$ws.Cells.Item(122, 14).Value = -14759.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(35, 8).Value = 173333.33
$ws.Cells.Item(35, 9).Value = 5000
$ws.Cells.Item(35, 11).Value = 5000
$ws.Cells.Item(35, 13).Value = -4664
$ws.Cells.Item(46, 8).Value = 6783.7856
$ws.Cells.Item(46, 10).Value = 7271.8184
$ws.Cells.Item(46, 12).Value = 7271.8184
$ws.Cells.Item(46, 14).Value = -7647.8184

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 1599.75
$ws.Cells.Item(2, 9).Value = 1599.75
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 1599.75
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -1487.75
$ws.Cells.Item(2, 14).ClearContents()
$ws.Cells.Item(32, 8).Value = 40176.668
$ws.Cells.Item(32, 9).Value = 8500
$ws.Cells.Item(32, 10).Value = 56015
$ws.Cells.Item(32, 11).Value = 8500
$ws.Cells.Item(32, 12).Value = 56015
$ws.Cells.Item(32, 13).Value = -8183
$ws.Cells.Item(32, 14).Value = -56649
$ws.Cells.Item(107, 8).Value = 503.53845
$ws.Cells.Item(107, 10).Value = 500
$ws.Cells.Item(107, 12).Value = 1500
$ws.Cells.Item(107, 14).Value = -5340
$ws.Cells.Item(132, 8).Value = 1873.5385
$ws.Cells.Item(132, 9).Value = 1869.2424
$ws.Cells.Item(132, 10).Value = 1897.1666
$ws.Cells.Item(132, 11).Value = 5607.7272
$ws.Cells.Item(132, 12).Value = 5691.4998
$ws.Cells.Item(132, 13).Value = -3077.7272
$ws.Cells.Item(132, 14).Value = -10751.4998
